$d = $word.ActiveDocument

# --- Create the three new character styles (matching styles.xml additions) ---
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Underline = 1
$gaNLinks.Font.Color = 8388608

# --- Update every "Orionin tähtikuvio havainnointijaksot ..." run: add a ---
# --- trailing period and apply the GaNStyle character style.            ---
$searchText = "Orionin tähtikuvio havainnointijaksot vuonna 2022: 16.-25.1., 14.-23.2., 14.-24.3"
$rng = $d.Content
while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.InsertAfter(".")
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

# --- Apply GaNParagraph to the "Osallistut maailmanlaajuiseen ..." run ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpiä näkyvissä olevia tähtiä keinona mitata valonsaastetta tietyssä paikassa. Paikallistamalla ja tarkkailemalla Orionin tähtikuvio miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot päivittyvät heti verkossa olevaan tietokantaan, ja näin saadaan käsitys siitä minkä verran taivaan tähdistä on missäkin nähtävissä.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Tämän oppaan kartat piirsi Jenik ..." run ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Tämän oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Style = "GaNLinks"
}
